# RGossF-HW10.xlsx update: rerun of simulation with two additional samples
# (Holden, Rizzie Spiral) inserted right after Spiral5, and a renamed sample
# (Thomas Hex -> Matthies Hex).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new data rows right after the "Spiral5" row (row 4) ---
# This shifts the existing rows 4:29 down to 6:31, preserving all of their
# values/format, and leaves two blank rows (4:5) ready for the new samples.
$ws.Rows("4:5").Insert()

# Re-apply the same direct formatting used by the rest of column A (bold,
# bordered, centered) to the two newly inserted A-cells so they match their
# neighbours exactly (Insert() alone leaves them with a blended style).
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. New sample "Holden" (row 4) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.504623523186097
$ws.Range("D4").Value = 0.8367063306086788
$ws.Range("E4").Value = 0.8367063306086788
$ws.Range("F4").Value = 0.9145542588308017
$ws.Range("G4").Value = 1.002269474040153
$ws.Range("H4").Value = 0.9053427232817619
$ws.Range("I4").Value = 0.8889496893244269
$ws.Range("J4").Value = 1.504623523186097
$ws.Range("K4").Value = 1.504623523186097
$ws.Range("L4").Value = 0.9145542588308017
$ws.Range("M4").Value = 0.8756302947197403
$ws.Range("N4").Value = 0.8756302947197403
$ws.Range("O4").Value = 0.8800700929213026
$ws.Range("P4").Value = 1.085294704208526
$ws.Range("Q4").Value = 1.085294704208526
$ws.Range("R4").Value = 1.190126908952919
$ws.Range("S4").Value = 1.190126908952919
$ws.Range("T4").Value = 1.008740999878653

# --- 3. New sample "Rizzie Spiral" (row 5) ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.570478024798651
$ws.Range("D5").Value = 3.792226568893989
$ws.Range("E5").Value = 3.792226568893989
$ws.Range("F5").Value = 0.06554019872590683
$ws.Range("G5").Value = 0.3317472814760187
$ws.Range("H5").Value = 0.001087735600771453
$ws.Range("I5").Value = 2.070066004530333
$ws.Range("J5").Value = 1.570478024798651
$ws.Range("K5").Value = 1.570478024798651
$ws.Range("L5").Value = 0.06554019872590683
$ws.Range("M5").Value = 1.928883383809948
$ws.Range("N5").Value = 1.928883383809948
$ws.Range("O5").Value = 1.97594425738341
$ws.Range("P5").Value = 1.809414930806182
$ws.Range("Q5").Value = 1.809414930806182
$ws.Range("R5").Value = 1.749680704304299
$ws.Range("S5").Value = 1.749680704304299
$ws.Range("T5").Value = 1.305190969004278

# --- 4. Rename sample "Thomas Hex" -> "Matthies Hex" ---
# After the row insert above, the row that used to be row 9 ("Thomas Hex")
# is now row 11.
$ws.Range("B11").Value = "Matthies Hex"

Write-Host "Edit applied: inserted Holden/Rizzie Spiral rows, renamed Thomas Hex to Matthies Hex."
